$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing credential row (vishal7singh03 / wefggew) into row 3
$ws.Range("A3").Value = "vishal7singh03"
$ws.Range("B3").Value = "wefggew"

# Add a new Instagram credential pair in row 4
$ws.Range("A4").Value = "r2442tt"
$ws.Range("B4").Value = "ggg42g2"
